$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B-column filenames for rows 53-75 (A and other columns follow the existing
# pattern: A = "n" + (row-1), C = True, D = no_meltpatch, E = negative)
$fileNames = @(
    "n52_IMG_3178HorFlip.jpeg",
    "n53_IMG_3178HorVertFlip.jpeg",
    "n54_IMG_3178VertFlip.jpeg",
    "n55_IMG_3179.jpeg",
    "n56_IMG_3179HorFlip.jpeg",
    "n57_IMG_3179HorVertFlip.jpeg",
    "n58_IMG_3179VertFlip.jpeg",
    "n59_IMG_3180.jpeg",
    "n60_IMG_3180HorFlip.jpeg",
    "n61_IMG_3180HorVertFlip.jpeg",
    "n62_IMG_3180VertFlip.jpeg",
    "n63_IMG_3070.jpeg",
    "n64_IMG_3070HorFlip.jpeg",
    "n65_IMG_3070HorVertFlip.jpeg",
    "n66_IMG_3070VertFlip.jpeg",
    "n67_IMG_3072.jpeg",
    "n68_IMG_3072HorFlip.jpeg",
    "n69_IMG_3072HorVertFlip.jpeg",
    "n70_IMG_3072VertFlip.jpeg",
    "n71_IMG_3073.jpeg",
    "n72_IMG_3073HorFlip.jpeg",
    "n73_IMG_3073HorVertFlip.jpeg",
    "n74_IMG_3073VertFlip.jpeg"
)

$startRow = 53
for ($i = 0; $i -lt $fileNames.Length; $i++) {
    $row = $startRow + $i
    $subjectNum = $row - 1

    $ws.Cells.Item($row, 1).Value = "n$subjectNum"
    $ws.Cells.Item($row, 2).Value = $fileNames[$i]
    # Prefix with an apostrophe so Excel stores the literal text "True"
    # instead of auto-converting it to a boolean value.
    $ws.Cells.Item($row, 3).Value = "'True"
    $ws.Cells.Item($row, 4).Value = "no_meltpatch"
    $ws.Cells.Item($row, 5).Value = "negative"
}

# Update the saved selection to match the new extent of data
$ws.Range("A40:E64").Select()
